# Fruta / hortaliza, semanal
# Update the weekly price rows for "Terminal La Palmera de La Serena - Tuna".
# Rows 10-12 get shuffled/re-dated values, and two more weekly rows (14 & 15)
# are appended, pushing the previous last row (old row 13) down to row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common/shared values for this block of rows.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100107
$producto   = "Otros"
$categoriaId = 100107011
$categoria   = "Tuna"
$variedad    = "Sin especificar"
$unidad      = "`$/caja 18 kilos"
$origen      = "Provincia de Limarí"
$kgUnidad    = 18

function Set-TunaRow {
    param(
        $Row,
        $Fecha,
        $Calidad,
        $Volumen,
        $PrecioMin,
        $PrecioMax,
        $PrecioProm,
        $PrecioKg
    )

    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $productoId
    $ws.Cells.Item($Row, 8).Value  = $producto
    $ws.Cells.Item($Row, 9).Value  = $categoriaId
    $ws.Cells.Item($Row, 10).Value = $categoria
    $ws.Cells.Item($Row, 11).Value = $variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

Set-TunaRow 10 44609 "Primera"  240 13000 14000 13500 750
Set-TunaRow 11 44609 "Segunda"  240 11000 12000 11500 639
Set-TunaRow 12 44294 "Especial" 200 14500 15000 14750 819
Set-TunaRow 13 44294 "Primera"  240 12500 13000 12750 708
Set-TunaRow 14 44294 "Segunda"  240 10500 11000 10750 597
Set-TunaRow 15 44595 "Primera"  200 15500 16000 15750 875
